# "zoom in and out of image" (id 3) moves from the Active sheet to the
# Inactive sheet, as the new first task row there ("v1: zoom in and out;
# fit to window").

$wb = $excel.ActiveWorkbook
$wsActive = $wb.Worksheets.Item("Active")
$wsInactive = $wb.Worksheets.Item("Inactive")

# Row 3 on "Active" holds id=3 "zoom in and out of image".
$srcRow = 3

$idVal       = $wsActive.Cells.Item($srcRow, 1).Value()
$titleVal    = $wsActive.Cells.Item($srcRow, 2).Value()
$statusVal   = $wsActive.Cells.Item($srcRow, 3).Value()
$categoryVal = $wsActive.Cells.Item($srcRow, 4).Value()
$createdVal  = $wsActive.Cells.Item($srcRow, 5).Value()

# Remove it from "Active" - remaining rows shift up.
$wsActive.Rows($srcRow).Delete()

# Insert a fresh row at the top of the "Inactive" data (row 2, right under
# the header) - existing rows shift down.
$wsInactive.Rows(2).Insert()

# Created/Done hold date-looking text ("8/9/2018"); force text formatting
# first so Excel doesn't reinterpret it as a date serial number.
$wsInactive.Range("E2:F2").NumberFormat = "@"

$wsInactive.Cells.Item(2, 1).Value = $idVal
$wsInactive.Cells.Item(2, 2).Value = $titleVal
$wsInactive.Cells.Item(2, 3).Value = $statusVal
$wsInactive.Cells.Item(2, 4).Value = $categoryVal
$wsInactive.Cells.Item(2, 5).Value = $createdVal
$wsInactive.Cells.Item(2, 6).Value = $createdVal

# Match the plain (non-header) row styling used by the rest of the table.
$wsInactive.Range("A2:F2").Style = "Normal"
